$d = $word.ActiveDocument

# Use constants for Find.Execute parameters (mirrors Word's wdReplaceAll etc.)
$wdReplaceNone = 0
$wdReplaceOne = 1
$wdReplaceAll = 2

function Replace-UniqueText($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, $wdReplaceOne)
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-UniqueText "Unraveling the Enigmatic Universe" "Unraveling the Secrets of Life: A Journey Into the Realm of Biology"

# ---------------------------------------------------------------------------
# 2. Author name
# ---------------------------------------------------------------------------
Replace-UniqueText " Amelia Anderson" " Sarah Michelle"

# ---------------------------------------------------------------------------
# 3. Email user part ("amelia" -> "info@sarahmichelle")
# ---------------------------------------------------------------------------
Replace-UniqueText "amelia" "info@sarahmichelle"

# ---------------------------------------------------------------------------
# 4. Email domain part: merges "anderson@celestialobservatory" + "." + "org"
#    runs into a single run with text "com" (the preceding "." run between
#    "amelia" and "anderson@..." stays untouched).
# ---------------------------------------------------------------------------
Replace-UniqueText "anderson@celestialobservatory.org" "com"

# ---------------------------------------------------------------------------
# 5. First body paragraph, first sentence (its own run; trailing "." run
#    stays separate/unchanged).
# ---------------------------------------------------------------------------
Replace-UniqueText "As we gaze upon the vast expanse of the cosmos, a symphony of celestial bodies unfolds before our eyes" "Biology, the study of life, takes us on an awe-inspiring journey into the intricate workings of organisms, both large and small"

# ---------------------------------------------------------------------------
# 6. Second sentence of first paragraph -- merges 4 runs (the original
#    second/third/fourth sentences incl. their periods) into a single run,
#    leaving the period run that follows (before the <w:br/>) untouched.
# ---------------------------------------------------------------------------
Replace-UniqueText " From the intricate ballet of planets orbiting stars to the grand cosmic dance of galaxies, the universe reveals its awe-inspiring beauty and complexity. Yet, this enchanting tapestry also harbors profound mysteries that challenge our understanding and ignite our curiosity. In this exploration, we embark on a journey to unravel the enigmas of the universe, delving into the captivating realms of dark matter, black holes, and the elusive concept of time" " From the delicate cells that make up our bodies to the awe-inspiring diversity of ecosystems that surround us, biology offers a captivating lens through which we can understand the fundamental principles that govern life on Earth"

# ---------------------------------------------------------------------------
# 7. "The cosmos holds..." sentence (run directly after the double <w:br/>).
# ---------------------------------------------------------------------------
Replace-UniqueText "The cosmos holds a tantalizing mystery within its vastness: dark matter" "In this captivating exploration, we unravel the mysteries of life's origins, delving into the evolutionary processes that have shaped the breathtaking diversity of organisms"

# ---------------------------------------------------------------------------
# 8. Next sentence (own run; trailing "." stays separate).
# ---------------------------------------------------------------------------
Replace-UniqueText " This invisible substance, believed to constitute over 80% of the universe, exerts a gravitational influence that shapes the structure and evolution of galaxies" " We probe the inner workings of cells, uncovering the intricate dance of organelles that sustain life's essential functions"

# ---------------------------------------------------------------------------
# 9. Merge the following two sentences (and the period between them) into a
#    single run; the trailing "." after them stays separate/unchanged.
# ---------------------------------------------------------------------------
Replace-UniqueText " Its enigmatic nature has puzzled scientists for decades, as it interacts with ordinary matter in ways we do not fully comprehend. Unraveling the secrets of dark matter promises to illuminate the very fabric of our universe, casting light on its formation and destiny" " The world of genetics beckons us to unravel the secrets of heredity, as we marvel at the complex interplay of genes and their role in shaping our traits"

# ---------------------------------------------------------------------------
# 10. "Another cosmic enigma..." sentence (run directly after the double
#     <w:br/>).
# ---------------------------------------------------------------------------
Replace-UniqueText "Another cosmic enigma that captivates the scientific community is the black hole" "As we delve deeper into the realm of biology, we uncover the extraordinary complexity of ecosystems"

# ---------------------------------------------------------------------------
# 11. Next sentence (own run; trailing "." stays separate).
# ---------------------------------------------------------------------------
Replace-UniqueText " These celestial behemoths, with their immense gravitational pull, warp space-time and create regions of intense density" " From microscopic organisms to towering trees, each species occupies a niche, contributing to the intricate web of life"

# ---------------------------------------------------------------------------
# 12. Delete the large dead block of sentences that follows (runs covering
#     "As matter approaches a black hole..." through "...This " along with
#     both <w:br/> elements and the <w:lastRenderedPageBreak/> run), using a
#     Range so the markup (breaks / page-break hint) is removed outright
#     rather than merged into a surviving run.
# ---------------------------------------------------------------------------
$anchorStart = $d.Content
$anchorStart.Find.Execute("As matter approaches a black hole") | Out-Null
$deleteFrom = $anchorStart.Start

$anchorEnd = $d.Range($deleteFrom, $d.Content.End)
$anchorEnd.Find.Execute("interplay between space, time, and gravity poses fundamental questions about the nature of reality and the fundamental laws that govern our universe") | Out-Null
$deleteTo = $anchorEnd.Start

$d.Range($deleteFrom, $deleteTo).Text = ""

# ---------------------------------------------------------------------------
# 13. Replace the remaining "interplay between..." text (now immediately
#     following " This ") with the new sentence; this also naturally drops
#     the <w:lastRenderedPageBreak/> hint since that run gets merged away.
# ---------------------------------------------------------------------------
Replace-UniqueText "interplay between space, time, and gravity poses fundamental questions about the nature of reality and the fundamental laws that govern our universe" " We investigate the intricate relationships between organisms, exploring the delicate balance that maintains stability and resilience in our natural world"

# ---------------------------------------------------------------------------
# 14. Summary paragraph sentences (periods between/after remain separate,
#     unchanged runs).
# ---------------------------------------------------------------------------
Replace-UniqueText "In our exploration of the universe's enigmas, we have delved into the mysteries of dark matter, black holes, and the nature of time" "Our journey into biology has unveiled the wonders of life, from the microscopic world of cells to the vast expanse of ecosystems"

Replace-UniqueText " The uncharted territories of these cosmic frontiers hold the key to unlocking profound insights into the very fabric of our existence" " We have explored the evolutionary tapestry that has woven together all living things, unraveled the mysteries of heredity, and marveled at the interconnectedness of organisms in diverse environments"

Replace-UniqueText " By unraveling these enigmas, we embark on a journey to comprehend the universe's grand design, pushing the boundaries of human knowledge and unveiling the secrets that lie hidden within the vast expanse of the cosmos" " Biology continues to challenge our understanding of life while offering awe-inspiring insights into the fundamental principles that govern the existence of all organisms"

# ---------------------------------------------------------------------------
# 15. Append a new, fully empty paragraph at the very end of the document
#     body (before the sectPr).
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.Text = "`r"
